$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (A1:G12) is sorted by the "number" column (A) and a new LeetCode
# entry -- "Two Sum II - Input Array Is Sorted" (#167) -- was solved and
# needs to be inserted in sorted order between #147 (row 9) and #206
# (old row 10). Insert a fresh row at 10 and shift everything below down.
$ws.Rows.Item(10).Insert() | Out-Null

# Fill in the new row's data.
$ws.Range("A10").Value = 167
$ws.Range("B10").Value = "Medium"
$ws.Range("C10").Value = "Two Sum II - Input Array Is Sorted"
$ws.Range("D10").Value = "http://rb.gy/psjwn9"
$ws.Range("E10").Value = "Pointers"
$ws.Range("F10").Value = "O(n)"
$ws.Range("G10").Value = "Use left and right pointers and the fact that the array is sorted"

# (Row.Insert() already carries the correct per-column styles into the new
# row, so no explicit style assignment is needed for A10/B10/C10/E10/F10/G10.)

# The row Insert() shifts the worksheet cells down but this emulation does
# not shift hyperlink anchors along with it, so the hyperlinks that used to
# live on D10/D11/D12 are still anchored there (now pointing at the wrong
# row). Re-home them one row lower, and give the new row's D10 cell the
# hyperlink for the newly added problem.
$ws.Hyperlinks.Add($ws.Range("D13"), "http://rb.gy/nrugfa") | Out-Null
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$D$12') { $h.Address = "http://rb.gy/bcqtel" }
    elseif ($addr -eq '$D$11') { $h.Address = "http://rb.gy/1nj72g" }
    elseif ($addr -eq '$D$10') { $h.Address = "http://rb.gy/psjwn9" }
}

# Typing a new row nudges the "next" empty-row selection down by one.
$ws.Range("C16").Select() | Out-Null
